# Update cryptocurrency list (price + 1h volume change) as scraped on
# Fri Jul 12 14:00:23 UTC 2024.
#
# Note: several "Price" values look like plain decimals (e.g. "0.999",
# "531.03") which Excel would otherwise silently auto-convert to numbers.
# To keep them as literal text (matching the source inlineStr cells), a
# leading apostrophe is used to force text entry, and the cell style is
# reset back to "Normal" right afterwards so no stray number-format style
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.933.22'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").Value = '3.122.46'
$ws.Range("E3").Value = '  -1.63%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''531.03'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.66%  '

$ws.Range("D6").Value = '''139.09'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("D7").Value = '''0.997'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.34%  '

$ws.Range("D8").Value = '3.122.85'
$ws.Range("E8").Value = '  -1.56%  '

$ws.Range("E9").Value = '  +3.84%  '

$ws.Range("E10").Value = '  +0.27%  '

$ws.Range("E11").Value = '  -2.64%  '

$ws.Range("D12").Value = '''0.410'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +2.33%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.654.80'
$ws.Range("E13").Value = '  -1.68%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '''0.137'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +1.72%  '

$ws.Range("D15").Value = '''25.57'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("D16").Value = '''0.0000163'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -2.56%  '

$ws.Range("D17").Value = '57.862.00'
$ws.Range("E17").Value = '  -1.17%  '

$ws.Range("D18").Value = '3.113.91'
$ws.Range("E18").Value = '  -1.76%  '

$ws.Range("D19").Value = '''5.96'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -3.53%  '

$ws.Range("D20").Value = '''12.63'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -2.35%  '

$ws.Range("D21").Value = '''8.00'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").Value = '''351.83'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -2.12%  '

$ws.Range("E23").Value = '  +0.20%  '

$ws.Range("D24").Value = '''69.04'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").Value = '''0.505'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -1.97%  '

$ws.Range("E26").Value = '  -2.09%  '

$ws.Range("D27").Value = '''0.996'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.53%  '

$ws.Range("D28").Value = '0.0₃0877'
$ws.Range("E28").Value = '  -8.53%  '

$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '''7.20'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -4.19%  '

$ws.Range("E31").Value = '  -1.81%  '

$ws.Range("D32").Value = '''6.03'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -7.37%  '

$ws.Range("D33").Value = '''21.26'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").Value = '''4.96'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("D35").Value = '''1.15'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -6.48%  '

$ws.Range("D36").Value = '''158.76'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +0.71%  '

$ws.Range("D37").Value = '''6.07'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -2.36%  '

$ws.Range("D38").Value = '''25.87'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("D39").Value = '''1.26'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -4.04%  '

$ws.Range("D40").Value = '''1.69'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +1.94%  '

$ws.Range("D41").Value = '''0.0669'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -1.23%  '

$ws.Range("D42").Value = '''4.05'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +0.29%  '

$ws.Range("D43").Value = '''0.697'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -1.66%  '

$ws.Range("D44").Value = '2.409.59'
$ws.Range("E44").Value = '  +1.71%  '

$ws.Range("D45").Value = '''36.95'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").Value = '3.158.32'
$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("D48").Value = '''0.0265'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -2.75%  '

$ws.Range("D49").Value = '''0.961'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -5.65%  '

$ws.Range("D50").Value = '''6.05'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -0.71%  '

$ws.Range("D51").Value = '''19.80'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -4.37%  '

